$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row data matching the target diff (rows 2-11)
$data = @(
    @{Row=2;  A=41390; B="Luana da Rosa";      C="Financeiro";       D="Viagem de negócios"; E=1; F=45082; G=5563.83},
    @{Row=3;  A=36311; B="Catarina Duarte";    C="Recursos Humanos"; D="Problemas pessoais"; E=3; F=45085; G=4477.47},
    @{Row=4;  A=1068;  B="Isabelly Moreira";   C="Jurídico";         D="Consulta médica";    E=2; F=45091; G=6154.62},
    @{Row=5;  A=39902; B="Davi Lucca Correia"; C="Engenharia";       D="Problemas pessoais"; E=1; F=45081; G=7379.41},
    @{Row=6;  A=45101; B="Luiz Felipe da Luz"; C="Jurídico";         D="Doença";             E=8; F=45106; G=10797.49},
    @{Row=7;  A=46812; B="Gabrielly Oliveira"; C="Engenharia";       D="Viagem de negócios"; E=6; F=45085; G=12057.58},
    @{Row=8;  A=20786; B="Nicolas Lopes";      C="Jurídico";         D="Doença";             E=5; F=45103; G=8819.48},
    @{Row=9;  A=67104; B="Benjamin Fernandes"; C="Jurídico";         D="Doença";             E=7; F=45080; G=3830.61},
    @{Row=10; A=6365;  B="Diogo da Costa";     C="Marketing";        D="Viagem de negócios"; E=8; F=45096; G=10476.64},
    @{Row=11; A=15453; B="Srta. Lara Moura";   C="P&D";              D="Outros";             E=4; F=45086; G=9897.83}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}

$wb.Save()
